$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark the numeric-looking (D/E column) cells as Text format so
# Excel stores the new values as strings (matching the source inlineStr cells)
# instead of coercing "244.91" / "-0.73%" into numbers.
$numFmtRange = $ws.Range("D2:E50")
$numFmtRange.NumberFormat = "@"

# Step 2: write the updated cell values
$ws.Range("D2").Value = "244.91"
$ws.Range("E2").Value = "-0.73%"
$ws.Range("D3").Value = "27.46"
$ws.Range("E3").Value = "5.10%"
$ws.Range("D4").Value = "5.114"
$ws.Range("E4").Value = "0.55%"
$ws.Range("D5").Value = "0.05685"
$ws.Range("E5").Value = "1.66%"
$ws.Range("D6").Value = "6.517"
$ws.Range("E6").Value = "0.69%"
$ws.Range("D7").Value = "0.8196"
$ws.Range("E7").Value = "0.73%"
$ws.Range("D8").Value = "0.8570"
$ws.Range("E8").Value = "1.36%"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "0.01016"
$ws.Range("E9").Value = "1,603.54%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1336"
$ws.Range("E10").Value = "0.23%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.06944"
$ws.Range("E11").Value = "-0.47%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.02877"
$ws.Range("E12").Value = "2.12%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09395"
$ws.Range("E13").Value = "0.09%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001511"
$ws.Range("E14").Value = "-0.27%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "0.04025"
$ws.Range("E15").Value = "-13.77%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006215"
$ws.Range("E16").Value = "0.11%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.511"
$ws.Range("E17").Value = "-2.68%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "3.009"
$ws.Range("E18").Value = "-0.36%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.319"
$ws.Range("E19").Value = "12.83%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3189"
$ws.Range("E20").Value = "2.47%"
$ws.Range("D21").Value = "0.03221"
$ws.Range("E21").Value = "0.53%"
$ws.Range("E22").Value = "-0.04%"
$ws.Range("D23").Value = "3.553"
$ws.Range("E23").Value = "-5.85%"
$ws.Range("E24").Value = "-0.13%"
$ws.Range("E25").Value = "-2.62%"
$ws.Range("D26").Value = "0.004479"
$ws.Range("E26").Value = "-1.87%"
$ws.Range("E27").Value = "22.83%"
$ws.Range("E28").Value = "-27.44%"
$ws.Range("D40").Value = "0.03717"
$ws.Range("E41").Value = "75.62%"
$ws.Range("D42").Value = "0.1059"
$ws.Range("E42").Value = "-22.53%"
$ws.Range("D43").Value = "0.002164"
$ws.Range("E43").Value = "-17.15%"
$ws.Range("D44").Value = "0.009711"
$ws.Range("E44").Value = "20.29%"
$ws.Range("D45").Value = "0.00005119"
$ws.Range("E45").Value = "-5.00%"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("E47").Value = "-30.36%"
$ws.Range("D48").Value = "0.002517"
$ws.Range("E48").Value = "3.95%"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("E50").Value = "-0.02%"

# Step 3: drop the temporary Text number format so the cells end up with
# the default (unstyled) formatting, same as the source file.
$numFmtRange.ClearFormats()
